$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert 2 new blank columns before column B (shifts B:BV -> D:BX)
$ws.Range("B:C").Insert()

# Step 2: set new C1 label
$ws.Range("C1").Value = "RASCHEL"

# Step 3: append trailing space to the RASCHEL2 section labels (now at BC1:BN1 and BP1)
$ws.Range("BC1").Value = "F9 "
$ws.Range("BD1").Value = "F12 "
$ws.Range("BE1").Value = "F14 "
$ws.Range("BF1").Value = "F18 "
$ws.Range("BG1").Value = "F9TB "
$ws.Range("BH1").Value = "F12TB "
$ws.Range("BI1").Value = "F14TB "
$ws.Range("BJ1").Value = "F18TB "
$ws.Range("BK1").Value = "F9TC "
$ws.Range("BL1").Value = "F12TC "
$ws.Range("BM1").Value = "F14TC "
$ws.Range("BN1").Value = "F18TC "
$ws.Range("BP1").Value = "T3975 "

# Step 4: remove the stray single-space cell (shift left)
$ws.Range("BS1").Delete(-4159)

# Step 5: add row 2 data
$ws.Range("A2").Value = 12
$ws.Range("F2").Value = 4
$ws.Range("H2").Value = 5
$ws.Range("J2").Value = 5
$ws.Range("T2").Value = 5
$ws.Range("V2").Value = 3
$ws.Range("AB2").Value = 5
$ws.Range("AF2").Value = 12
$ws.Range("AG2").Value = 10
$ws.Range("AH2").Value = 5

Write-Output "done"
